$d = $word.ActiveDocument

$replacements = @(
    @("73÷2=", "47÷2="),
    @("45÷2=", "22÷6="),
    @("66÷6=", "38÷6="),
    @("87÷7=", "48÷8="),
    @("58÷5=", "72÷7="),
    @("54÷6=", "33÷2="),
    @("35÷5=", "56÷9="),
    @("62÷3=", "51÷9="),
    @("35÷8=", "12÷7="),
    @("70÷2=", "61÷8="),
    @("96÷4=", "36÷9="),
    @("75÷8=", "84÷4="),
    @("49÷6=", "99÷2="),
    @("75÷2=", "78÷5="),
    @("82÷9=", "16÷7="),
    @("53÷4=", "65÷6="),
    @("35÷6=", "30÷9="),
    @("91÷8=", "95÷3="),
    @("45÷3=", "71÷8="),
    @("23÷5=", "95÷7="),
    @("60÷5=", "79÷6="),
    @("29÷9=", "42÷7="),
    @("34÷9=", "39÷4="),
    @("19÷9=", "98÷8="),
    @("64÷9=", "10÷4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
